$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.394.87'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.910.23'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").Value = '''246.55'
$ws.Range("E5").Value = '  +2.60%  '
$ws.Range("D6").Value = '''0.667'
$ws.Range("E6").Value = '  +6.59%  '
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").Value = '''41.75'
$ws.Range("E8").Value = '  -2.55%  '
$ws.Range("E9").Value = '  +4.29%  '
$ws.Range("D10").Value = '''53.29'
$ws.Range("E10").Value = '  +13.49%  '
$ws.Range("D11").Value = '''0.0719'
$ws.Range("E11").Value = '  +2.78%  '
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("D13").Value = '2.189.15'
$ws.Range("E13").Value = '  +1.70%  '
$ws.Range("E14").Value = '  +6.12%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.702'
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.917.55'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("D18").Value = '35.405.76'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '''72.14'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").Value = '0.0₃0821'
$ws.Range("E20").Value = '  +2.22%  '
$ws.Range("D21").Value = '''241.62'
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").Value = '''12.49'
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").Value = '''4.83'
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").Value = '''2.39'
$ws.Range("E25").Value = '  +26.26%  '
$ws.Range("D26").Value = '''2.28'
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("D27").Value = '''171.15'
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").Value = '''8.47'
$ws.Range("E28").Value = '  +2.70%  '
$ws.Range("D29").Value = '''18.39'
$ws.Range("E29").Value = '  +3.04%  '
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("D31").Value = '4.156.66'
$ws.Range("E31").Value = '  +21.75%  '
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("D33").Value = '''0.953'
$ws.Range("E33").Value = '  +15.03%  '
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").Value = '''1.73'
$ws.Range("E37").Value = '  -4.82%  '
$ws.Range("D38").Value = '''2.04'
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").Value = '''1.34'
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.0653'
$ws.Range("E41").Value = '  +7.96%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0208'
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").Value = '''16.29'
$ws.Range("E43").Value = '  +6.81%  '
$ws.Range("D44").Value = '''89.97'
$ws.Range("E44").Value = '  -1.21%  '
$ws.Range("D45").Value = '1.340.89'
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("D46").Value = '''48.90'
$ws.Range("E46").Value = '  +39.11%  '
$ws.Range("E47").Value = '  +1.68%  '
$ws.Range("E48").Value = '  +2.21%  '
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("D50").Value = '''6.51'
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").Value = '2.096.84'
$ws.Range("E51").Value = '  +1.58%  '
